# Criando funcao que cria backup do arquivo que alimenta o pbi
# Update the numeric values in the "Valores" sheet (A1:L16) to reflect
# the refreshed report data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,4,1,5,0,2,0,0,3,0,7,0),
    @(0,0,0,0,1,1,0,3,0,5,1,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,1,0,0,0,0,0,0,1,0),
    @(0,3,0,0,0,0,0,0,0,0,0,0),
    @(1,5,1,4,0,2,0,0,4,0,6,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,1,0,1,0,0),
    @(0,0,0,0,2,1,0,1,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,1,0,0,0,0,1,1,1),
    @(2,12,2,10,4,6,0,5,7,7,16,2)
)

for ($r = 1; $r -le 16; $r++) {
    $rowValues = $data[$r - 1]
    for ($c = 1; $c -le 12; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
